$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "QUANTIDADE"

$values = @(1091,945,3951,12775,32188,63885,105711,152809,192543,219326,226837,218271,199207,173008,145976,121068,99501,81181,66329,55233,46223,39286,33230,28758,24597,21651,19212,16937,14965,13448,12012,10527,9498,8148,7195,6372,5588,4633,3986,3022,2382,1719,1176,569,224,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
